# Weekly update for "Hortaliza, Vega Modelo de Temuco - Camote":
# a new week's record is inserted at row 99 (pushing the existing history
# down by one row, so the oldest record that used to be the last one in
# the table now becomes the new last row), and the new top record gets
# this week's date while keeping last week's other figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 99:163 down to 100:164, inserting a new (initially blank,
# but format-inheriting) row 99.
$ws.Rows.Item(99).Insert()

# Populate the new row 99 with this week's record.
$ws.Range("A99").Value = 10
$ws.Range("B99").Value = "Vega Modelo de Temuco"
$ws.Range("C99").Value = "La Araucanía"
$ws.Range("D99").Value = 44957
$ws.Range("E99").Value = 9
$ws.Range("F99").Value = 100114002
$ws.Range("G99").Value = "Camote"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 20
$ws.Range("K99").Value = 24000
$ws.Range("L99").Value = 24000
$ws.Range("M99").Value = 24000
$ws.Range("N99").Value = "$/malla 20 kilos"
$ws.Range("O99").Value = "Perú"
$ws.Range("P99").Value = 1200
$ws.Range("Q99").Value = 20
$ws.Range("R99").Value = "Hortaliza"
